$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '41.574.63'
$ws.Cells.Item(2, 4).ClearFormats()

$ws.Cells.Item(2, 5).Value = '  -1.37%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.165.34'
$ws.Cells.Item(3, 4).ClearFormats()

$ws.Cells.Item(3, 5).Value = '  -2.76%  '

$ws.Cells.Item(4, 5).Value = '  -0.12%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '237.76'
$ws.Cells.Item(5, 4).ClearFormats()

$ws.Cells.Item(5, 5).Value = '  -2.15%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.600'
$ws.Cells.Item(6, 4).ClearFormats()

$ws.Cells.Item(6, 5).Value = '  -4.25%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '71.85'
$ws.Cells.Item(7, 4).ClearFormats()

$ws.Cells.Item(7, 5).Value = '  -3.21%  '

$ws.Cells.Item(8, 5).Value = '  -0.19%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.577'
$ws.Cells.Item(9, 4).ClearFormats()

$ws.Cells.Item(9, 5).Value = '  -4.42%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '39.62'
$ws.Cells.Item(10, 4).ClearFormats()

$ws.Cells.Item(10, 5).Value = '  -7.27%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0906'
$ws.Cells.Item(11, 4).ClearFormats()

$ws.Cells.Item(11, 5).Value = '  -5.63%  '

$ws.Cells.Item(12, 5).Value = '  -3.87%  '

$ws.Cells.Item(13, 5).Value = '  -2.65%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.69'
$ws.Cells.Item(14, 4).ClearFormats()

$ws.Cells.Item(14, 5).Value = '  -3.99%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.486.49'
$ws.Cells.Item(15, 4).ClearFormats()

$ws.Cells.Item(15, 5).Value = '  -2.89%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.26'
$ws.Cells.Item(16, 4).ClearFormats()

$ws.Cells.Item(16, 5).Value = '  -0.29%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.141.32'
$ws.Cells.Item(17, 4).ClearFormats()

$ws.Cells.Item(17, 5).Value = '  -3.33%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.775'
$ws.Cells.Item(18, 4).ClearFormats()

$ws.Cells.Item(18, 5).Value = '  -7.30%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '41.416.45'
$ws.Cells.Item(19, 4).ClearFormats()

$ws.Cells.Item(19, 5).Value = '  -1.40%  '

$ws.Cells.Item(20, 5).Value = '  -2.82%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '69.95'
$ws.Cells.Item(21, 4).ClearFormats()

$ws.Cells.Item(21, 5).Value = '  -3.97%  '

$ws.Cells.Item(22, 5).Value = '  -7.12%  '

$ws.Cells.Item(23, 5).Value = '  -10.19%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '225.91'
$ws.Cells.Item(24, 4).ClearFormats()

$ws.Cells.Item(24, 5).Value = '  -1.89%  '

$ws.Cells.Item(25, 5).Value = '  -3.97%  '

$ws.Cells.Item(26, 5).Value = '  -0.25%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.72'
$ws.Cells.Item(27, 4).ClearFormats()

$ws.Cells.Item(27, 5).Value = '  -5.79%  '

$ws.Cells.Item(28, 5).Value = '  -9.81%  '

$ws.Cells.Item(29, 5).Value = '  -3.95%  '

$ws.Cells.Item(30, 5).Value = '  -1.73%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '171.09'
$ws.Cells.Item(31, 4).ClearFormats()

$ws.Cells.Item(31, 5).Value = '  +2.45%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '19.79'
$ws.Cells.Item(32, 4).ClearFormats()

$ws.Cells.Item(32, 5).Value = '  -3.93%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '32.83'
$ws.Cells.Item(33, 4).ClearFormats()

$ws.Cells.Item(33, 5).Value = '  +9.31%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0769'
$ws.Cells.Item(34, 4).ClearFormats()

$ws.Cells.Item(34, 5).Value = '  -4.09%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.36'
$ws.Cells.Item(35, 4).ClearFormats()

$ws.Cells.Item(35, 5).Value = '  -4.87%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.119'
$ws.Cells.Item(36, 4).ClearFormats()

$ws.Cells.Item(36, 5).Value = '  -4.33%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.28'
$ws.Cells.Item(37, 4).ClearFormats()

$ws.Cells.Item(37, 5).Value = '  -1.27%  '

$ws.Cells.Item(38, 5).Value = '  -7.28%  '

$ws.Cells.Item(39, 5).Value = '  -0.30%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '12.10'
$ws.Cells.Item(40, 4).ClearFormats()

$ws.Cells.Item(40, 5).Value = '  -8.48%  '

$ws.Cells.Item(41, 5).Value = '  -1.94%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.36'
$ws.Cells.Item(42, 4).ClearFormats()

$ws.Cells.Item(42, 5).Value = '  -5.98%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '58.93'
$ws.Cells.Item(43, 4).ClearFormats()

$ws.Cells.Item(43, 5).Value = '  -9.41%  '

$ws.Cells.Item(44, 5).Value = '  -2.91%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.188'
$ws.Cells.Item(45, 4).ClearFormats()

$ws.Cells.Item(45, 5).Value = '  -5.39%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0966'
$ws.Cells.Item(46, 4).ClearFormats()

$ws.Cells.Item(46, 5).Value = '  -3.66%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '96.94'
$ws.Cells.Item(47, 4).ClearFormats()

$ws.Cells.Item(47, 5).Value = '  -7.11%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.08'
$ws.Cells.Item(48, 4).ClearFormats()

$ws.Cells.Item(48, 5).Value = '  -3.92%  '

$ws.Cells.Item(49, 5).Value = '  -5.00%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.18'
$ws.Cells.Item(50, 4).ClearFormats()

$ws.Cells.Item(50, 5).Value = '  -7.40%  '

$ws.Cells.Item(51, 5).Value = '  -2.44%  '
